# Generate Report for Handoff
# Update "Latest Handoff Datetime" (column D) for rows that were just
# (re-)handed off: rows 4, 6, 7, 8, 9, 10 on both the zh-cn and de-de
# status sheets. Rows 2/3 are already "Handed back: in sync with en-US"
# and row 5 is "In Translation" - those keep their existing timestamps.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$zhcnTimestamp = "2016-03-10 10:06:39"
$ws_zhcn.Range("D4").Value = $zhcnTimestamp
$ws_zhcn.Range("D6").Value = $zhcnTimestamp
$ws_zhcn.Range("D7").Value = $zhcnTimestamp
$ws_zhcn.Range("D8").Value = $zhcnTimestamp
$ws_zhcn.Range("D9").Value = $zhcnTimestamp
$ws_zhcn.Range("D10").Value = $zhcnTimestamp

$ws_dede = $wb.Worksheets.Item("de-de")
$dedeTimestamp = "2016-03-10 10:06:51"
$ws_dede.Range("D4").Value = $dedeTimestamp
$ws_dede.Range("D6").Value = $dedeTimestamp
$ws_dede.Range("D7").Value = $dedeTimestamp
$ws_dede.Range("D8").Value = $dedeTimestamp
$ws_dede.Range("D9").Value = $dedeTimestamp
$ws_dede.Range("D10").Value = $dedeTimestamp
